# Refresh the scraped "Price" (D) / "Volume(1h)" (E) columns with a new
# crawl snapshot, matching the GitHub Actions "Updated symbol list" commit.
#
# The source cells are plain text (scraper writes strings, not numbers),
# so a bare `.Value = "303.95"` would let Excel's type-inference coerce it
# to a number. Prefixing with a literal apostrophe forces text entry (the
# same trick you'd use typing into the grid), then we drop the style back
# to "Normal" so the cell doesn't keep a quote-prefixed number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "303.95";       "E2"  = "5.19%"
    "D3"  = "35.07";        "E3"  = "12.94%"
    "D4"  = "5.177";        "E4"  = "4.44%"
    "D5"  = "0.07823";      "E5"  = "6.32%"
    "D6"  = "2.296";        "E6"  = "-1.67%"
                            "E7"  = "4.37%"
    "D8"  = "3.987";        "E8"  = "6.92%"
    "D9"  = "0.9234";       "E9"  = "1.41%"
    "D10" = "0.1006";       "E10" = "8.80%"
    "D11" = "0.1836";       "E11" = "8.34%"
    "D12" = "0.08488";      "E12" = "2.47%"
    "D13" = "0.03379";      "E13" = "8.44%"
    "D14" = "0.09915";      "E14" = "-0.55%"
    "D15" = "0.001487";     "E15" = "-0.37%"
    "D16" = "0.04651";      "E16" = "3.09%"
    "D17" = "0.005799";     "E17" = "1.68%"
    "D18" = "3.483";        "E18" = "-0.28%"
    "D19" = "2.128";        "E19" = "1.57%"
    "D20" = "0.3441";       "E20" = "3.46%"
    "D21" = "0.1326";       "E21" = "3.35%"
    "D22" = "4.584";        "E22" = "9.72%"
    "D23" = "0.2396";       "E23" = "14.10%"
    "D24" = "0.001223";     "E24" = "0.97%"
    "D25" = "0.004328";     "E25" = "3.59%"
    "D26" = "0.0001301";    "E26" = "0.12%"
    "D27" = "0.0003402";    "E27" = "0.28%"
    "D39" = "0.01738";      "E39" = "10.32%"
    "D40" = "0.04744";      "E40" = "6.32%"
                            "E41" = "4.45%"
    "D42" = "0.1410";       "E42" = "5.89%"
    "D43" = "0.007264";     "E43" = "-26.72%"
    "D44" = "0.002292";     "E44" = "1.90%"
    "D45" = "0.01007";      "E45" = "15.06%"
    "D46" = "0.00006033";   "E46" = "-1.17%"
    "D47" = "0.00000000751";"E47" = "0.12%"
    "D48" = "5.797";        "E48" = "125.97%"
    "D49" = "0.002692";     "E49" = "34.64%"
    "D50" = "0.00002102";   "E50" = "0.12%"
    "D51" = "0.0002001";    "E51" = "0.12%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
